$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "admin"
$ws.Range("C2").Value = 0

$ws.Range("F4").Select()
